$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1998.8823
$ws.Range("I28").Value = 2026.36
$ws.Range("J28").Value = 1922.5555
$ws.Range("K28").Value = 2026.36
$ws.Range("L28").Value = 1922.5555
$ws.Range("M28").Value = -1541.36
$ws.Range("N28").Value = -2892.5555

$ws.Range("H33").Value = 459.875
$ws.Range("I33").Value = 367
$ws.Range("J33").Value = 515.6
$ws.Range("K33").Value = 367
$ws.Range("L33").Value = 515.6
$ws.Range("M33").Value = -138
$ws.Range("N33").Value = -973.6

$ws.Range("H41").Value = 1298
$ws.Range("I41").Value = 1118.2727
$ws.Range("K41").Value = 1118.2727
$ws.Range("M41").Value = -678.2727

$ws.Range("H113").Value = 1987.5454
$ws.Range("I113").Value = 1987
$ws.Range("J113").Value = 1988.5
$ws.Range("K113").Value = 1987
$ws.Range("L113").Value = 1988.5
$ws.Range("M113").Value = 1267
$ws.Range("N113").Value = -8496.5

$ws.Range("H116").Value = 2627.0908
$ws.Range("I116").Value = 2100
$ws.Range("J116").Value = 2928.2856
$ws.Range("K116").Value = 2100
$ws.Range("L116").Value = 2928.2856
$ws.Range("M116").Value = 1342
$ws.Range("N116").Value = -9812.285599999999

$ws.Range("H132").Value = 11503022
$ws.Range("I132").Value = 15158866
$ws.Range("K132").Value = 45476598
$ws.Range("M132").Value = -45474068

$ws.Range("H137").Value = 1473.5
$ws.Range("I137").Value = 1179.5
$ws.Range("J137").Value = 1996.1666
$ws.Range("K137").Value = 3538.5
$ws.Range("L137").Value = 5988.4998
$ws.Range("M137").Value = -988.5
$ws.Range("N137").Value = -11088.4998

$ws.Range("H138").Value = 2367.9333
$ws.Range("I138").Value = 5133.3335
$ws.Range("J138").Value = 2252.7083
$ws.Range("K138").Value = 15400.0005
$ws.Range("L138").Value = 6758.124899999999
$ws.Range("M138").Value = -10260.0005
$ws.Range("N138").Value = -17038.1249


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 945.3182
$ws.Range("I2").Value = 677.6111
$ws.Range("K2").Value = 677.6111
$ws.Range("M2").Value = -564.6111

$ws.Range("H32").Value = 6464.6743
$ws.Range("I32").Value = 6531.318
$ws.Range("K32").Value = 6531.318
$ws.Range("M32").Value = -6244.318

$ws.Range("H61").Value = 38462756
$ws.Range("I61").Value = 43479264
$ws.Range("J61").Value = 2871.3333
$ws.Range("K61").Value = 43479264
$ws.Range("L61").Value = 2871.3333
$ws.Range("M61").Value = -43479052
$ws.Range("N61").Value = -3295.3333

$ws.Range("H102").Value = 10418039
$ws.Range("I102").Value = 12821548
$ws.Range("J102").Value = 2833.3333
$ws.Range("K102").Value = 12821548
$ws.Range("L102").Value = 2833.3333
$ws.Range("M102").Value = -12819926
$ws.Range("N102").Value = -6077.3333

$ws.Range("H116").Value = 945.3182
$ws.Range("I116").Value = 677.6111
$ws.Range("K116").Value = 677.6111
$ws.Range("M116").Value = 1616.3889

$ws.Range("H122").Value = 1973.2069
$ws.Range("I122").Value = 1928.36
$ws.Range("K122").Value = 5785.08
$ws.Range("M122").Value = -3335.08

$ws.Range("H136").Value = 38462756
$ws.Range("I136").Value = 43479264
$ws.Range("J136").Value = 2871.3333
$ws.Range("K136").Value = 130437792
$ws.Range("L136").Value = 8613.999899999999
$ws.Range("M136").Value = -130435242
$ws.Range("N136").Value = -13713.9999


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 945.3182
$ws.Range("I3").Value = 677.6111
$ws.Range("K3").Value = 677.6111
$ws.Range("M3").Value = -563.6111

$ws.Range("H107").Value = 922.4231
$ws.Range("I107").Value = 658.5
$ws.Range("K107").Value = 658.5
$ws.Range("M107").Value = 1261.5

$ws.Range("H135").Value = 35363.332
$ws.Range("J135").Value = 35363.332
$ws.Range("L135").Value = 35363.332
$ws.Range("N135").Value = -45503.332


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 62501160
$ws.Range("I16").Value = 83334520
$ws.Range("K16").Value = 83334520
$ws.Range("M16").Value = -83334233

$ws.Range("H31").Value = 1195.3677
$ws.Range("I31").Value = 1087.3387
$ws.Range("K31").Value = 1087.3387
$ws.Range("M31").Value = -792.3387

$ws.Range("H34").Value = 1195.3677
$ws.Range("I34").Value = 1087.3387
$ws.Range("K34").Value = 1087.3387
$ws.Range("M34").Value = -885.3387

$ws.Range("H113").Value = 62501160
$ws.Range("I113").Value = 83334520
$ws.Range("K113").Value = 83334520
$ws.Range("M113").Value = -83332350

$ws.Range("H135").Value = 35466.363
$ws.Range("J135").Value = 35466.363
$ws.Range("L135").Value = 35466.363
$ws.Range("N135").Value = -45606.363


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 192.3
$ws.Range("I2").Value = 203
$ws.Range("J2").Value = 187.71428
$ws.Range("K2").Value = 1218
$ws.Range("L2").Value = 1126.28568
$ws.Range("M2").Value = -1105
$ws.Range("N2").Value = -1352.28568

$ws.Range("H34").Value = 5883678.5
$ws.Range("I34").Value = 393.14285
$ws.Range("J34").Value = 10001978
$ws.Range("K34").Value = 1179.42855
$ws.Range("L34").Value = 30005934
$ws.Range("M34").Value = -1095.42855
$ws.Range("N34").Value = -30006102

$ws.Range("H131").Value = 20439532
$ws.Range("J131").Value = 40307.344
$ws.Range("L131").Value = 120922.032
$ws.Range("N131").Value = -131002.032

$ws.Range("H139").Value = 1838.6666
$ws.Range("I139").Value = 1796.6875
$ws.Range("J139").Value = 1973
$ws.Range("K139").Value = 5390.0625
$ws.Range("L139").Value = 5919
$ws.Range("M139").Value = -250.0625
$ws.Range("N139").Value = -16199


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 14893.75
$ws.Range("J92").Value = 14893.75
$ws.Range("L92").Value = 14893.75
$ws.Range("N92").Value = -18637.75

$ws.Range("H102").Value = 1469.3823
$ws.Range("I102").Value = 1429
$ws.Range("K102").Value = 1429
$ws.Range("M102").Value = 193

$ws.Range("H126").Value = 2050.1667
$ws.Range("I126").Value = 1784.0834
$ws.Range("J126").Value = 2582.3333
$ws.Range("K126").Value = 5352.2502
$ws.Range("L126").Value = 7746.999899999999
$ws.Range("M126").Value = -2882.2502
$ws.Range("N126").Value = -12686.9999


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 693.06665
$ws.Range("I22").Value = 818.5
$ws.Range("J22").Value = 609.44446
$ws.Range("K22").Value = 818.5
$ws.Range("L22").Value = 609.44446
$ws.Range("M22").Value = -523.5
$ws.Range("N22").Value = -1199.44446

$ws.Range("H27").Value = 693.06665
$ws.Range("I27").Value = 818.5
$ws.Range("J27").Value = 609.44446
$ws.Range("K27").Value = 818.5
$ws.Range("L27").Value = 609.44446
$ws.Range("M27").Value = -711.5
$ws.Range("N27").Value = -823.44446

$ws.Range("H68").Value = 1231
$ws.Range("I68").Value = 1227.3684
$ws.Range("K68").Value = 1227.3684
$ws.Range("M68").Value = -478.3684000000001

$ws.Range("H71").Value = 1231
$ws.Range("I71").Value = 1227.3684
$ws.Range("K71").Value = 6136.842000000001
$ws.Range("M71").Value = -2392.842000000001

$ws.Range("H122").Value = 27779446
$ws.Range("I122").Value = 35715856
$ws.Range("K122").Value = 107147568
$ws.Range("M122").Value = -107145118


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 446.25
$ws.Range("I100").Value = 464
$ws.Range("J100").Value = 393
$ws.Range("K100").Value = 928
$ws.Range("L100").Value = 786
$ws.Range("M100").Value = -387
$ws.Range("N100").Value = -1868

$ws.Range("H122").Value = 29071220
$ws.Range("I122").Value = 32896174
$ws.Range("J122").Value = 1558
$ws.Range("K122").Value = 98688522
$ws.Range("L122").Value = 4674
$ws.Range("M122").Value = -98686072
$ws.Range("N122").Value = -9574

